$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column R by copying column Q (values + formatting), then
# overwrite the 2021-year figures. Using a full-column copy+insert keeps
# the exact same style (cellXfs) indices as column Q for every row,
# matching what Excel does when a user duplicates an existing column.
$ws.Columns("Q").Copy()
$ws.Columns("R").Insert(-4161)

# Update the new column's header to the new year
$ws.Range("R4").Value = 2021

# Update the data values for 2021 (rows where real figures are available).
# Rows left untouched (3, 11, 17, 18, 19, 29-34) keep the "-" placeholder
# / blank formatting copied from column Q.
$ws.Range("R5").Value = 109
$ws.Range("R6").Value = 74
$ws.Range("R7").Value = 35
$ws.Range("R8").Value = 36
$ws.Range("R9").Value = 35
$ws.Range("R10").Value = 1
$ws.Range("R12").Value = 8
$ws.Range("R13").Value = 7
$ws.Range("R14").Value = 12
$ws.Range("R15").Value = 7
$ws.Range("R16").Value = 5
$ws.Range("R20").Value = 17
$ws.Range("R21").Value = 8
$ws.Range("R22").Value = 9
$ws.Range("R23").Value = 9
$ws.Range("R24").Value = 7
$ws.Range("R25").Value = 2
$ws.Range("R26").Value = 20
$ws.Range("R27").Value = 9
$ws.Range("R28").Value = 11

# Match the author's final selection/active cell
$ws.Range("R1").Select()
